$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content changes (headers) ---
$ws.Range("A2").Value = "Seedling species"
$ws.Range("B1").Value = "Model AICc value"

# --- Merge the header cell across B1:D1 ---
$ws.Range("B1:D1").Merge()

# --- Header row 1 (B1:D1): bold TNR12 already set; add bottom border + center align ---
$hdr1 = $ws.Range("B1:D1")
$hdr1.Borders.Item(9).LineStyle = 1
$hdr1.Borders.Item(9).Weight = 2
$hdr1.HorizontalAlignment = -4108

# --- Row 2, column A header: add right border (bottom border already present) ---
$a2 = $ws.Range("A2")
$a2.Borders.Item(10).LineStyle = 1
$a2.Borders.Item(10).Weight = 2

# --- Shaded rows 3-6 (grey fill matching theme Background2 darker 10%) ---
$greyColor = 13619152

$ws.Range("A3:A6").Interior.Color = $greyColor
$ws.Range("B3:D6").Interior.Color = $greyColor

# A3: right + top border (top border marks the header/body divider together with B1:D1's/A2's bottom)
$a3 = $ws.Range("A3")
$a3.Borders.Item(10).LineStyle = 1
$a3.Borders.Item(10).Weight = 2
$a3.Borders.Item(8).LineStyle = 1
$a3.Borders.Item(8).Weight = 2

# A4:A6: right border only
$a4_6 = $ws.Range("A4:A6")
$a4_6.Borders.Item(10).LineStyle = 1
$a4_6.Borders.Item(10).Weight = 2

# --- Row 7 (Aglaia mariannensis) - no fill, right border only on column A ---
$a7 = $ws.Range("A7")
$a7.Borders.Item(10).LineStyle = 1
$a7.Borders.Item(10).Weight = 2

# --- Row 8 (Ochrosia oppositifolia) - no fill, right border (bottom already present) ---
$a8 = $ws.Range("A8")
$a8.Borders.Item(10).LineStyle = 1
$a8.Borders.Item(10).Weight = 2

# --- Selection state ---
$ws.Range("A1:D8").Select()
